{"js": "const replacements = [\n  [\"2026-01-07 Wednesday\", \"2026-01-08 Thursday\"],\n  [\"60\u00d719=1140\", \"84\u00d740=3360\"],\n  [\"58\u00d773=4234\", \"23\u00d788=2024\"],\n  [\"26\u00d722=572\", \"43\u00d782=3526\"],\n  [\"27\u00d736=972\", \"50\u00d717=850\"],\n  [\"49\u00d787=4263\", \"38\u00d763=2394\"],\n  [\"23\u00d783=1909\", \"65\u00d733=2145\"],\n  [\"51\u00d723=1173\", \"71\u00d768=4828\"],\n  [\"79\u00d778=6162\", \"73\u00d725=1825\"],\n  [\"37\u00d727=999\", \"98\u00d748=4704\"],\n  [\"47\u00d728=1316\", \"89\u00d731=2759\"],\n  [\"91\u00d782=7462\", \"46\u00d788=4048\"],\n  [\"46\u00d727=1242\", \"87\u00d771=6177\"],\n  [\"25\u00d799=2475\", \"75\u00d778=5850\"],\n  [\"23\u00d754=1242\", \"43\u00d758=2494\"],\n  [\"46\u00d728=1288\", \"72\u00d792=6624\"],\n  [\"19\u00d769=1311\", \"55\u00d762=3410\"],\n  [\"81\u00d734=2754\", \"96\u00d777=7392\"],\n  [\"79\u00d792=7268\", \"92\u00d764=5888\"],\n  [\"13\u00d743=559\", \"47\u00d783=3901\"],\n  [\"87\u00d726=2262\", \"22\u00d774=1628\"],\n  [\"20\u00d798=1960\", \"59\u00d780=4720\"],\n  [\"68\u00d727=1836\", \"54\u00d772=3888\"],\n  [\"90\u00d712=1080\", \"41\u00d742=1722\"],\n  [\"49\u00d752=2548\", \"88\u00d734=2992\"],\n  [\"50\u00d772=3600\", \"67\u00d760=4020\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-07 Wednesday\", \"2026-01-08 Thursday\"),\n    @(\"60\u00d719=1140\", \"84\u00d740=3360\"),\n    @(\"58\u00d773=4234\", \"23\u00d788=2024\"),\n    @(\"26\u00d722=572\", \"43\u00d782=3526\"),\n    @(\"27\u00d736=972\", \"50\u00d717=850\"),\n    @(\"49\u00d787=4263\", \"38\u00d763=2394\"),\n    @(\"23\u00d783=1909\", \"65\u00d733=2145\"),\n    @(\"51\u00d723=1173\", \"71\u00d768=4828\"),\n    @(\"79\u00d778=6162\", \"73\u00d725=1825\"),\n    @(\"37\u00d727=999\", \"98\u00d748=4704\"),\n    @(\"47\u00d728=1316\", \"89\u00d731=2759\"),\n    @(\"91\u00d782=7462\", \"46\u00d788=4048\"),\n    @(\"46\u00d727=1242\", \"87\u00d771=6177\"),\n    @(\"25\u00d799=2475\", \"75\u00d778=5850\"),\n    @(\"23\u00d754=1242\", \"43\u00d758=2494\"),\n    @(\"46\u00d728=1288\", \"72\u00d792=6624\"),\n    @(\"19\u00d769=1311\", \"55\u00d762=3410\"),\n    @(\"81\u00d734=2754\", \"96\u00d777=7392\"),\n    @(\"79\u00d792=7268\", \"92\u00d764=5888\"),\n    @(\"13\u00d743=559\", \"47\u00d783=3901\"),\n    @(\"87\u00d726=2262\", \"22\u00d774=1628\"),\n    @(\"20\u00d798=1960\", \"59\u00d780=4720\"),\n    @(\"68\u00d727=1836\", \"54\u00d772=3888\"),\n    @(\"90\u00d712=1080\", \"41\u00d742=1722\"),\n    @(\"49\u00d752=2548\", \"88\u00d734=2992\"),\n    @(\"50\u00d772=3600\", \"67\u00d760=4020\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}"}
